$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1494.6154
$ws.Range("I19").Value = 863.93335
$ws.Range("J19").Value = 2354.6365
$ws.Range("K19").Value = 863.93335
$ws.Range("L19").Value = 2354.6365
$ws.Range("M19").Value = -688.93335
$ws.Range("N19").Value = -2704.6365

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3000.25
$ws.Range("I76").Value = 3000.25
$ws.Range("K76").Value = 3000.25
$ws.Range("M76").Value = -2685.25

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3000.25
$ws.Range("I79").Value = 3000.25
$ws.Range("K79").Value = 3000.25
$ws.Range("M79").Value = -1908.25

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2464.7058
$ws.Range("I80").Value = 1610.4445
$ws.Range("K80").Value = 4831.333500000001
$ws.Range("M80").Value = -3833.333500000001

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 2464.7058
$ws.Range("I83").Value = 1610.4445
$ws.Range("K83").Value = 14494.0005
$ws.Range("M83").Value = -9502.0005

# ALC row 88
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 7343.5
$ws.Range("I88").Value = 7295.3335
$ws.Range("J88").Value = 7364.143
$ws.Range("K88").Value = 7295.3335
$ws.Range("L88").Value = 7364.143
$ws.Range("M88").Value = -6889.3335
$ws.Range("N88").Value = -8176.143

# ALC row 91
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H91").Value = 7343.5
$ws.Range("I91").Value = 7295.3335
$ws.Range("J91").Value = 7364.143
$ws.Range("K91").Value = 7295.3335
$ws.Range("L91").Value = 7364.143
$ws.Range("M91").Value = -5891.3335
$ws.Range("N91").Value = -10172.143

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2679.5518
$ws.Range("I132").Value = 2748.04
$ws.Range("J132").Value = 2251.5
$ws.Range("K132").Value = 8244.119999999999
$ws.Range("L132").Value = 6754.5
$ws.Range("M132").Value = -5714.119999999999
$ws.Range("N132").Value = -11814.5

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1324.1833
$ws.Range("I61").Value = 1202.9556
$ws.Range("J61").Value = 1687.8667
$ws.Range("K61").Value = 1202.9556
$ws.Range("L61").Value = 1687.8667
$ws.Range("M61").Value = -990.9556
$ws.Range("N61").Value = -2111.8667

# ARM row 88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2495.5557
$ws.Range("I88").Value = 2320
$ws.Range("K88").Value = 2320
$ws.Range("M88").Value = -1914

# ARM row 91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2495.5557
$ws.Range("I91").Value = 2320
$ws.Range("K91").Value = 2320
$ws.Range("M91").Value = -916

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2107
$ws.Range("I122").Value = 1768.5
$ws.Range("J122").Value = 2558.3333
$ws.Range("K122").Value = 5305.5
$ws.Range("L122").Value = 7674.999899999999
$ws.Range("M122").Value = -2855.5
$ws.Range("N122").Value = -12574.9999

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1945.5294
$ws.Range("I132").Value = 1688.1818
$ws.Range("J132").Value = 2417.3333
$ws.Range("K132").Value = 5064.5454
$ws.Range("L132").Value = 7251.999899999999
$ws.Range("M132").Value = -2534.5454
$ws.Range("N132").Value = -12311.9999

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1324.1833
$ws.Range("I136").Value = 1202.9556
$ws.Range("J136").Value = 1687.8667
$ws.Range("K136").Value = 3608.8668
$ws.Range("L136").Value = 5063.6001
$ws.Range("M136").Value = -1058.8668
$ws.Range("N136").Value = -10163.6001

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2184.5386
$ws.Range("I20").Value = 2066.6667
$ws.Range("J20").Value = 2285.5715
$ws.Range("K20").Value = 2066.6667
$ws.Range("L20").Value = 2285.5715
$ws.Range("M20").Value = -1819.6667
$ws.Range("N20").Value = -2779.5715

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1720.8889
$ws.Range("I86").Value = 1568.85
$ws.Range("J86").Value = 2155.2856
$ws.Range("K86").Value = 1568.85
$ws.Range("L86").Value = 2155.2856
$ws.Range("M86").Value = -445.8499999999999
$ws.Range("N86").Value = -4401.2856

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1720.8889
$ws.Range("I89").Value = 1568.85
$ws.Range("J89").Value = 2155.2856
$ws.Range("K89").Value = 7844.25
$ws.Range("L89").Value = 10776.428
$ws.Range("M89").Value = -2228.25
$ws.Range("N89").Value = -22008.428

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1342.4348
$ws.Range("I99").Value = 1168.8
$ws.Range("K99").Value = 1168.8
$ws.Range("M99").Value = 329.2

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 550875.9
$ws.Range("I134").Value = 802578.0600000001
$ws.Range("J134").Value = 3697.1738
$ws.Range("K134").Value = 2407734.18
$ws.Range("L134").Value = 11091.5214
$ws.Range("M134").Value = -2405199.18
$ws.Range("N134").Value = -16161.5214

# CRP row 44
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 59500
$ws.Range("I44").Value = 59000
$ws.Range("J44").Value = 60000
$ws.Range("K44").Value = 59000
$ws.Range("L44").Value = 60000
$ws.Range("M44").Value = -58558
$ws.Range("N44").Value = -60884

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1147.2972
$ws.Range("I58").Value = 1074.5333
$ws.Range("J58").Value = 1459.1428
$ws.Range("K58").Value = 1074.5333
$ws.Range("L58").Value = 1459.1428
$ws.Range("M58").Value = -871.5333000000001
$ws.Range("N58").Value = -1865.1428

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1766596.1
$ws.Range("I132").Value = 3110.7058
$ws.Range("J132").Value = 9261409
$ws.Range("K132").Value = 9332.117400000001
$ws.Range("L132").Value = 27784227
$ws.Range("M132").Value = -6802.117400000001
$ws.Range("N132").Value = -27789287

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1917.7963
$ws.Range("I134").Value = 2039.45
$ws.Range("J134").Value = 1570.2142
$ws.Range("K134").Value = 6118.35
$ws.Range("L134").Value = 4710.642599999999
$ws.Range("M134").Value = -3583.35
$ws.Range("N134").Value = -9780.642599999999

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1147.2972
$ws.Range("I136").Value = 1074.5333
$ws.Range("J136").Value = 1459.1428
$ws.Range("K136").Value = 3223.5999
$ws.Range("L136").Value = 4377.428400000001
$ws.Range("M136").Value = -673.5999000000002
$ws.Range("N136").Value = -9477.428400000001

# GSM row 48
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 11350
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 11350
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 11350
$ws.Range("N48").Value = -12320
$ws.Range("M48").ClearContents()

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4615.372
$ws.Range("I70").Value = 4459.457
$ws.Range("J70").Value = 5297.5
$ws.Range("K70").Value = 4459.457
$ws.Range("L70").Value = 5297.5
$ws.Range("M70").Value = -4189.457
$ws.Range("N70").Value = -5837.5

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 4615.372
$ws.Range("I73").Value = 4459.457
$ws.Range("J73").Value = 5297.5
$ws.Range("K73").Value = 4459.457
$ws.Range("L73").Value = 5297.5
$ws.Range("M73").Value = -3523.457
$ws.Range("N73").Value = -7169.5

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2000
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -1002
$ws.Range("N80").ClearContents()

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2000
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -5008
$ws.Range("N83").ClearContents()

# GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1024.1
$ws.Range("I97").Value = 1053.3334
$ws.Range("J97").Value = 980.25
$ws.Range("K97").Value = 1053.3334
$ws.Range("L97").Value = 980.25
$ws.Range("M97").Value = -557.3334
$ws.Range("N97").Value = -1972.25

# GSM row 113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1494.4286
$ws.Range("I113").Value = 1338.3636
$ws.Range("J113").Value = 2066.6667
$ws.Range("K113").Value = 1338.3636
$ws.Range("L113").Value = 2066.6667
$ws.Range("M113").Value = 831.6364000000001
$ws.Range("N113").Value = -6406.6667

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 256888.75
$ws.Range("I122").Value = 505000
$ws.Range("K122").Value = 1515000
$ws.Range("M122").Value = -1512550

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2634283.5
$ws.Range("I132").Value = 2656.577
$ws.Range("J132").Value = 8336141.5
$ws.Range("K132").Value = 7969.731000000001
$ws.Range("L132").Value = 25008424.5
$ws.Range("M132").Value = -5439.731000000001
$ws.Range("N132").Value = -25013484.5

# LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 907.6087
$ws.Range("I46").Value = 1086.5714
$ws.Range("K46").Value = 1086.5714
$ws.Range("M46").Value = -898.5714

# LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 3000
$ws.Range("I48").Value = 3000
$ws.Range("K48").Value = 3000
$ws.Range("M48").Value = -2339

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5046.488
$ws.Range("I132").Value = 6001.1113
$ws.Range("K132").Value = 18003.3339
$ws.Range("M132").Value = -15473.3339

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 64200.895
$ws.Range("I122").Value = 13654.267
$ws.Range("J122").Value = 253750.75
$ws.Range("K122").Value = 40962.801
$ws.Range("L122").Value = 761252.25
$ws.Range("M122").Value = -38512.801
$ws.Range("N122").Value = -766152.25
